# GTY_YR_FIN.xlsx update: add a new "most-recent period" column.
#
# The source data in this sheet is laid out as several stacked tables
# (Income Statement / Balance Sheet / Cash Flow Statement), each with a
# "Period Ending" header row followed by line items. Every table uses the
# same column layout: C = label, D..J = 7 periods (existing), K = spare
# (blank) column. This edit adds one more (more recent) reporting period
# by inserting a new column before D, which pushes the existing D:K data
# to E:L, and then fills the new column D with the latest period's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 102

# 1) Insert a new blank column before column D. This shifts D:K -> E:L
#    (values, formulas and formatting all move together).
$ws.Columns("D:D").Insert()

# 2) The newly inserted column D currently inherits formatting from column
#    C (label column). Re-stripe it from column E (which holds what used
#    to be column D) so every row's new D cell gets the same number
#    format/font/alignment as the rest of that row's data (date style for
#    the header rows, the right-aligned number style for data rows, and
#    plain/blank for the spacer rows). Skip the lone table-title rows
#    (5, 6, 37, 79) that only have a label in column A/B and no data
#    columns at all.
$ws.Range("E7:E35").Copy() | Out-Null
$ws.Range("D7:D35").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("E38:E77").Copy() | Out-Null
$ws.Range("D38:D77").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("E80:E102").Copy() | Out-Null
$ws.Range("D80:D102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 3) Populate the new column D with the latest period's figures.
#    Header rows: new "Period Ending" date (31-Dec-2018 => serial 43465).
$ws.Range("D7").Value  = 43465
$ws.Range("D38").Value = 43465
$ws.Range("D80").Value = 43465

# -- Income Statement --------------------------------------------------
$ws.Range("D8").Value  = 136100
$ws.Range("D9").Value  = 28400
$ws.Range("D10").Value = 107700
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 4900
$ws.Range("D15").Value = 23600
$ws.Range("D17").Value = 72000
$ws.Range("D18").Value = 64100
$ws.Range("D20").Value = 6700
$ws.Range("D21").Value = 94400
$ws.Range("D22").Value = 22300
$ws.Range("D23").Value = 48400
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 48400
$ws.Range("D27").Value = 47700
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -700
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -6700
$ws.Range("D33").Value = 47000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 47000

# -- Balance Sheet ------------------------------------------------------
$ws.Range("D41").Value = 46900
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 40700
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 57900
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 119400
$ws.Range("D48").Value = 892400
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 1900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1159200
$ws.Range("D57").Value = 62100
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 14500
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 441600
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 578000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -57400
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 581200
$ws.Range("D77").Value = 0

# -- Cash Flow Statement -------------------------------------------------
$ws.Range("D81").Value = 47000
$ws.Range("D83").Value = 23600
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 63300
$ws.Range("D91").Value = -84400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -75900
$ws.Range("D96").Value = -50500
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 40500
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 27900

# Rows 11, 16, 19, 39, 40, 55, 56, 67, 82, 90, 95 are spacer rows that were
# (and remain) blank in every period column, so column D needs no value
# there - the format-only paste in step 2 already matches the rest of the
# row.

# Disable multi-threaded calculation (matches the workbook's
# concurrentCalc="0" calculation setting).
$excel.MultiThreadedCalculation.Enabled = $false
